$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "20.429.67"
$ws.Range("E2").Value = "  -7.19%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.441.12"
$ws.Range("E3").Value = "  -7.22%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "278.05"
$ws.Range("E6").Value = "  -3.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3755"
$ws.Range("E7").Value = "  -4.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3065"
$ws.Range("E8").Value = "  -3.81%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "40.46"
$ws.Range("E9").Value = "  -8.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.011"
$ws.Range("E10").Value = "  -4.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06552"
$ws.Range("E11").Value = "  -8.31%  "
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.366"
$ws.Range("E13").Value = "  -4.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.28"
$ws.Range("E14").Value = "  -6.62%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.139"
$ws.Range("E15").Value = "  -6.90%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.441.94"
$ws.Range("E16").Value = "  -7.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001010"
$ws.Range("E17").Value = "  -7.65%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.05872"
$ws.Range("E18").Value = "  -10.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "76.26"
$ws.Range("E19").Value = "  -7.68%  "
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.725"
$ws.Range("E21").Value = "  -6.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.40"
$ws.Range("E22").Value = "  -5.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.87"
$ws.Range("E23").Value = "  -1.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.310"
$ws.Range("E24").Value = "  -3.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "20.431.96"
$ws.Range("E25").Value = "  -7.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "143.01"
$ws.Range("E26").Value = "  -2.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.214"
$ws.Range("E27").Value = "  -5.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.01"
$ws.Range("E28").Value = "  -7.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.604.63"
$ws.Range("E29").Value = "  -7.16%  "
$ws.Range("E30").Value = "  -6.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.932"
$ws.Range("E31").Value = "  -19.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9052"
$ws.Range("E32").Value = "  -6.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.428"
$ws.Range("E33").Value = "  -6.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07743"
$ws.Range("E34").Value = "  -6.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.301"
$ws.Range("E35").Value = "  -8.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.002"
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.80"
$ws.Range("E37").Value = "  +1.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05646"
$ws.Range("E38").Value = "  -5.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.143"
$ws.Range("E39").Value = "  -4.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.719"
$ws.Range("E40").Value = "  -6.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02046"
$ws.Range("E41").Value = "  -8.36%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1916"
$ws.Range("E42").Value = "  -5.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.357"
$ws.Range("E43").Value = "  -14.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.586"
$ws.Range("E44").Value = "  -4.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5325"
$ws.Range("E45").Value = "  -7.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.15"
$ws.Range("E46").Value = "  -5.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5155"
$ws.Range("E47").Value = "  -6.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "111.83"
$ws.Range("E48").Value = "  -3.93%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.792"
$ws.Range("E49").Value = "  -3.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.053"
$ws.Range("E50").Value = "  -6.26%  "
$ws.Range("B51").Value = "PaxDollar"
$ws.Range("C51").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.003"
$ws.Range("E51").Value = "  +0.17%  "

Write-Host "Applied cryptos update"